$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(6, 8).Value = 3000
$ws.Cells.Item(6, 9).Value = 0
$ws.Cells.Item(6, 10).Value = 3000
$ws.Cells.Item(6, 11).Value = 0
$ws.Cells.Item(6, 12).Value = 9000
$ws.Cells.Item(6, 13).ClearContents()
$ws.Cells.Item(6, 14).Value = -9224
$ws.Cells.Item(17, 8).Value = 1942.5652
$ws.Cells.Item(17, 10).Value = 1942.5652
$ws.Cells.Item(17, 12).Value = 5827.6956
$ws.Cells.Item(17, 14).Value = -6163.6956
$ws.Cells.Item(40, 8).Value = 16679579
$ws.Cells.Item(40, 10).Value = 35724692
$ws.Cells.Item(40, 12).Value = 35724692
$ws.Cells.Item(40, 14).Value = -35725042
$ws.Cells.Item(43, 8).Value = 2113764.2
$ws.Cells.Item(43, 9).Value = 5629205
$ws.Cells.Item(43, 10).Value = 4499.8
$ws.Cells.Item(43, 11).Value = 5629205
$ws.Cells.Item(43, 12).Value = 4499.8
$ws.Cells.Item(43, 13).Value = -5629136
$ws.Cells.Item(43, 14).Value = -4637.8
$ws.Cells.Item(58, 8).Value = 255
$ws.Cells.Item(58, 9).Value = 255
$ws.Cells.Item(58, 11).Value = 765
$ws.Cells.Item(58, 13).Value = -615
$ws.Cells.Item(64, 8).Value = 18940808
$ws.Cells.Item(64, 9).Value = 6627687.5
$ws.Cells.Item(64, 10).Value = 31253928
$ws.Cells.Item(64, 11).Value = 6627687.5
$ws.Cells.Item(64, 12).Value = 31253928
$ws.Cells.Item(64, 13).Value = -6627439.5
$ws.Cells.Item(64, 14).Value = -31254424
$ws.Cells.Item(67, 8).Value = 18940808
$ws.Cells.Item(67, 9).Value = 6627687.5
$ws.Cells.Item(67, 10).Value = 31253928
$ws.Cells.Item(67, 11).Value = 6627687.5
$ws.Cells.Item(67, 12).Value = 31253928
$ws.Cells.Item(67, 13).Value = -6626829.5
$ws.Cells.Item(67, 14).Value = -31255644
$ws.Cells.Item(74, 8).Value = 4666.6665
$ws.Cells.Item(74, 9).Value = 4500
$ws.Cells.Item(74, 11).Value = 4500
$ws.Cells.Item(74, 13).Value = -3564
$ws.Cells.Item(77, 8).Value = 4666.6665
$ws.Cells.Item(77, 9).Value = 4500
$ws.Cells.Item(77, 11).Value = 22500
$ws.Cells.Item(77, 13).Value = -17820
$ws.Cells.Item(82, 8).Value = 5185.364
$ws.Cells.Item(82, 9).Value = 1173.1666
$ws.Cells.Item(82, 11).Value = 3519.4998
$ws.Cells.Item(82, 13).Value = -3113.4998
$ws.Cells.Item(85, 8).Value = 5185.364
$ws.Cells.Item(85, 9).Value = 1173.1666
$ws.Cells.Item(85, 11).Value = 3519.4998
$ws.Cells.Item(85, 13).Value = -2115.4998
$ws.Cells.Item(86, 8).Value = 2195276.2
$ws.Cells.Item(86, 9).Value = 3598183
$ws.Cells.Item(86, 11).Value = 3598183
$ws.Cells.Item(86, 13).Value = -3597060
$ws.Cells.Item(89, 8).Value = 2195276.2
$ws.Cells.Item(89, 9).Value = 3598183
$ws.Cells.Item(89, 11).Value = 17990915
$ws.Cells.Item(89, 13).Value = -17985299
$ws.Cells.Item(99, 8).Value = 5316.4443
$ws.Cells.Item(99, 9).Value = 2836.6
$ws.Cells.Item(99, 10).Value = 8416.25
$ws.Cells.Item(99, 11).Value = 8509.8
$ws.Cells.Item(99, 12).Value = 25248.75
$ws.Cells.Item(99, 13).Value = -7011.799999999999
$ws.Cells.Item(99, 14).Value = -28244.75
$ws.Cells.Item(100, 8).Value = 1447.25
$ws.Cells.Item(100, 9).Value = 1179.8334
$ws.Cells.Item(100, 11).Value = 1179.8334
$ws.Cells.Item(100, 13).Value = -638.8334
$ws.Cells.Item(101, 8).Value = 659.8571
$ws.Cells.Item(101, 9).Value = 323.8
$ws.Cells.Item(101, 10).Value = 1500
$ws.Cells.Item(101, 11).Value = 971.4000000000001
$ws.Cells.Item(101, 12).Value = 4500
$ws.Cells.Item(101, 13).Value = 650.5999999999999
$ws.Cells.Item(101, 14).Value = -7744
$ws.Cells.Item(104, 8).Value = 596.75
$ws.Cells.Item(104, 9).Value = 596.75
$ws.Cells.Item(104, 10).Value = 0
$ws.Cells.Item(104, 11).Value = 1790.25
$ws.Cells.Item(104, 12).Value = 0
$ws.Cells.Item(104, 13).Value = -43.25
$ws.Cells.Item(104, 14).ClearContents()
$ws.Cells.Item(115, 8).Value = 583
$ws.Cells.Item(115, 9).Value = 499.5
$ws.Cells.Item(115, 10).Value = 750
$ws.Cells.Item(115, 11).Value = 1498.5
$ws.Cells.Item(115, 12).Value = 2250
$ws.Cells.Item(115, 13).Value = 68.5
$ws.Cells.Item(115, 14).Value = -5384
$ws.Cells.Item(118, 8).Value = 406.4
$ws.Cells.Item(118, 9).Value = 406.4
$ws.Cells.Item(118, 11).Value = 1219.2
$ws.Cells.Item(118, 13).Value = 437.8000000000002
$ws.Cells.Item(121, 8).Value = 4594.25
$ws.Cells.Item(121, 10).Value = 4594.25
$ws.Cells.Item(121, 12).Value = 13782.75
$ws.Cells.Item(121, 14).Value = -17276.75
$ws.Cells.Item(127, 8).Value = 580.44446
$ws.Cells.Item(127, 9).Value = 580.44446
$ws.Cells.Item(127, 11).Value = 1741.33338
$ws.Cells.Item(127, 13).Value = 3218.66662
$ws.Cells.Item(132, 8).Value = 13502.81
$ws.Cells.Item(132, 9).Value = 5848.2188
$ws.Cells.Item(132, 10).Value = 21404.322
$ws.Cells.Item(132, 11).Value = 17544.6564
$ws.Cells.Item(132, 12).Value = 64212.966
$ws.Cells.Item(132, 13).Value = -15014.6564
$ws.Cells.Item(132, 14).Value = -69272.966
$ws.Cells.Item(135, 8).Value = 2756.0645
$ws.Cells.Item(135, 9).Value = 925.2381
$ws.Cells.Item(135, 11).Value = 8327.1429
$ws.Cells.Item(135, 13).Value = -5792.142900000001
$ws.Cells.Item(138, 8).Value = 2372.818
$ws.Cells.Item(138, 9).Value = 1817.1904
$ws.Cells.Item(138, 10).Value = 2632.111
$ws.Cells.Item(138, 11).Value = 5451.5712
$ws.Cells.Item(138, 12).Value = 7896.333
$ws.Cells.Item(138, 13).Value = -311.5712000000003
$ws.Cells.Item(138, 14).Value = -18176.333

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 13030.861
$ws.Cells.Item(32, 9).Value = 13521.685
$ws.Cells.Item(32, 10).Value = 11165.733
$ws.Cells.Item(32, 11).Value = 13521.685
$ws.Cells.Item(32, 12).Value = 11165.733
$ws.Cells.Item(32, 13).Value = -13234.685
$ws.Cells.Item(32, 14).Value = -11739.733
$ws.Cells.Item(45, 8).Value = 5058.154
$ws.Cells.Item(45, 9).Value = 4217.5
$ws.Cells.Item(45, 10).Value = 5778.7144
$ws.Cells.Item(45, 11).Value = 4217.5
$ws.Cells.Item(45, 12).Value = 5778.7144
$ws.Cells.Item(45, 13).Value = -3840.5
$ws.Cells.Item(45, 14).Value = -6532.7144
$ws.Cells.Item(61, 8).Value = 6511.8374
$ws.Cells.Item(61, 9).Value = 6982.027
$ws.Cells.Item(61, 10).Value = 3612.3333
$ws.Cells.Item(61, 11).Value = 6982.027
$ws.Cells.Item(61, 12).Value = 3612.3333
$ws.Cells.Item(61, 13).Value = -6770.027
$ws.Cells.Item(61, 14).Value = -4036.3333
$ws.Cells.Item(63, 8).Value = 2201
$ws.Cells.Item(63, 9).Value = 1999.5
$ws.Cells.Item(63, 10).Value = 2402.5
$ws.Cells.Item(63, 11).Value = 1999.5
$ws.Cells.Item(63, 12).Value = 2402.5
$ws.Cells.Item(63, 13).Value = -1313.5
$ws.Cells.Item(63, 14).Value = -3774.5
$ws.Cells.Item(66, 8).Value = 2201
$ws.Cells.Item(66, 9).Value = 1999.5
$ws.Cells.Item(66, 10).Value = 2402.5
$ws.Cells.Item(66, 11).Value = 9997.5
$ws.Cells.Item(66, 12).Value = 12012.5
$ws.Cells.Item(66, 13).Value = -6565.5
$ws.Cells.Item(66, 14).Value = -18876.5
$ws.Cells.Item(74, 8).Value = 6945709
$ws.Cells.Item(74, 9).Value = 12501002
$ws.Cells.Item(74, 10).Value = 1592.75
$ws.Cells.Item(74, 11).Value = 12501002
$ws.Cells.Item(74, 12).Value = 1592.75
$ws.Cells.Item(74, 13).Value = -12500128
$ws.Cells.Item(74, 14).Value = -3340.75
$ws.Cells.Item(77, 8).Value = 6945709
$ws.Cells.Item(77, 9).Value = 12501002
$ws.Cells.Item(77, 10).Value = 1592.75
$ws.Cells.Item(77, 11).Value = 62505010
$ws.Cells.Item(77, 12).Value = 7963.75
$ws.Cells.Item(77, 13).Value = -62500642
$ws.Cells.Item(77, 14).Value = -16699.75
$ws.Cells.Item(97, 8).Value = 430.125
$ws.Cells.Item(97, 9).Value = 435.63635
$ws.Cells.Item(97, 10).Value = 369.5
$ws.Cells.Item(97, 11).Value = 435.63635
$ws.Cells.Item(97, 12).Value = 369.5
$ws.Cells.Item(97, 13).Value = 60.36365000000001
$ws.Cells.Item(97, 14).Value = -1361.5
$ws.Cells.Item(102, 8).Value = 549841.75
$ws.Cells.Item(102, 9).Value = 762394.75
$ws.Cells.Item(102, 11).Value = 762394.75
$ws.Cells.Item(102, 13).Value = -760772.75
$ws.Cells.Item(110, 8).Value = 2271149.2
$ws.Cells.Item(110, 9).Value = 2915906.2
$ws.Cells.Item(110, 11).Value = 2915906.2
$ws.Cells.Item(110, 13).Value = -2913861.2
$ws.Cells.Item(132, 8).Value = 5164.6587
$ws.Cells.Item(132, 9).Value = 2553
$ws.Cells.Item(132, 11).Value = 7659
$ws.Cells.Item(132, 13).Value = -5129
$ws.Cells.Item(136, 8).Value = 6511.8374
$ws.Cells.Item(136, 9).Value = 6982.027
$ws.Cells.Item(136, 10).Value = 3612.3333
$ws.Cells.Item(136, 11).Value = 20946.081
$ws.Cells.Item(136, 12).Value = 10836.9999
$ws.Cells.Item(136, 13).Value = -18396.081
$ws.Cells.Item(136, 14).Value = -15936.9999

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(22, 8).Value = 800
$ws.Cells.Item(22, 9).Value = 800
$ws.Cells.Item(22, 11).Value = 800
$ws.Cells.Item(22, 13).Value = -627
$ws.Cells.Item(94, 8).Value = 857116
$ws.Cells.Item(94, 9).Value = 1142421.6
$ws.Cells.Item(94, 11).Value = 1142421.6
$ws.Cells.Item(94, 13).Value = -1141970.6
$ws.Cells.Item(97, 8).Value = 15499.333
$ws.Cells.Item(97, 9).Value = 12599.2
$ws.Cells.Item(97, 11).Value = 12599.2
$ws.Cells.Item(97, 13).Value = -11608.2
$ws.Cells.Item(99, 8).Value = 1226821.2
$ws.Cells.Item(99, 9).Value = 1737196.9
$ws.Cells.Item(99, 11).Value = 1737196.9
$ws.Cells.Item(99, 13).Value = -1735698.9
$ws.Cells.Item(105, 8).Value = 2269.348
$ws.Cells.Item(105, 9).Value = 2045.7778
$ws.Cells.Item(105, 11).Value = 2045.7778
$ws.Cells.Item(105, 13).Value = -298.7778000000001
$ws.Cells.Item(107, 8).Value = 7199.8
$ws.Cells.Item(107, 9).Value = 7333
$ws.Cells.Item(107, 11).Value = 7333
$ws.Cells.Item(107, 13).Value = -5413

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 6849.523
$ws.Cells.Item(31, 9).Value = 2220.5
$ws.Cells.Item(31, 10).Value = 7878.1943
$ws.Cells.Item(31, 11).Value = 2220.5
$ws.Cells.Item(31, 12).Value = 7878.1943
$ws.Cells.Item(31, 13).Value = -1925.5
$ws.Cells.Item(31, 14).Value = -8468.1943
$ws.Cells.Item(34, 8).Value = 6849.523
$ws.Cells.Item(34, 9).Value = 2220.5
$ws.Cells.Item(34, 10).Value = 7878.1943
$ws.Cells.Item(34, 11).Value = 2220.5
$ws.Cells.Item(34, 12).Value = 7878.1943
$ws.Cells.Item(34, 13).Value = -2018.5
$ws.Cells.Item(34, 14).Value = -8282.1943
$ws.Cells.Item(62, 8).Value = 31234.455
$ws.Cells.Item(62, 9).Value = 3161.3333
$ws.Cells.Item(62, 11).Value = 3161.3333
$ws.Cells.Item(62, 13).Value = -2537.3333
$ws.Cells.Item(64, 8).Value = 0
$ws.Cells.Item(64, 10).Value = 0
$ws.Cells.Item(64, 12).Value = 0
$ws.Cells.Item(64, 14).ClearContents()
$ws.Cells.Item(65, 8).Value = 31234.455
$ws.Cells.Item(65, 9).Value = 3161.3333
$ws.Cells.Item(65, 11).Value = 15806.6665
$ws.Cells.Item(65, 13).Value = -12686.6665
$ws.Cells.Item(67, 8).Value = 0
$ws.Cells.Item(67, 10).Value = 0
$ws.Cells.Item(67, 12).Value = 0
$ws.Cells.Item(67, 14).ClearContents()
$ws.Cells.Item(94, 8).Value = 4363.143
$ws.Cells.Item(94, 9).Value = 0
$ws.Cells.Item(94, 10).Value = 4363.143
$ws.Cells.Item(94, 11).Value = 0
$ws.Cells.Item(94, 12).Value = 4363.143
$ws.Cells.Item(94, 13).ClearContents()
$ws.Cells.Item(94, 14).Value = -5265.143
$ws.Cells.Item(97, 8).Value = 35000
$ws.Cells.Item(97, 10).Value = 35000
$ws.Cells.Item(97, 12).Value = 35000
$ws.Cells.Item(97, 14).Value = -36982
$ws.Cells.Item(99, 8).Value = 13606.056
$ws.Cells.Item(99, 9).Value = 19944.285
$ws.Cells.Item(99, 11).Value = 19944.285
$ws.Cells.Item(99, 13).Value = -18446.285
$ws.Cells.Item(109, 8).Value = 77999.8
$ws.Cells.Item(109, 10).Value = 77999.8
$ws.Cells.Item(109, 12).Value = 77999.8
$ws.Cells.Item(109, 14).Value = -80079.8
$ws.Cells.Item(120, 8).Value = 44023.168
$ws.Cells.Item(120, 10).Value = 44023.168
$ws.Cells.Item(120, 12).Value = 44023.168
$ws.Cells.Item(120, 14).Value = -51281.168
$ws.Cells.Item(126, 8).Value = 13606.056
$ws.Cells.Item(126, 9).Value = 19944.285
$ws.Cells.Item(126, 11).Value = 59832.855
$ws.Cells.Item(126, 13).Value = -57362.855
$ws.Cells.Item(132, 8).Value = 13347715
$ws.Cells.Item(132, 9).Value = 13899703
$ws.Cells.Item(132, 10).Value = 100000
$ws.Cells.Item(132, 11).Value = 41699109
$ws.Cells.Item(132, 12).Value = 300000
$ws.Cells.Item(132, 13).Value = -41696579
$ws.Cells.Item(132, 14).Value = -305060
$ws.Cells.Item(134, 8).Value = 1702.0952
$ws.Cells.Item(134, 9).Value = 1375.4324
$ws.Cells.Item(134, 10).Value = 4119.4
$ws.Cells.Item(134, 11).Value = 4126.2972
$ws.Cells.Item(134, 12).Value = 12358.2
$ws.Cells.Item(134, 13).Value = -1591.2972
$ws.Cells.Item(134, 14).Value = -17428.2

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(7, 8).Value = 224.63637
$ws.Cells.Item(7, 9).Value = 224.14285
$ws.Cells.Item(7, 10).Value = 225.5
$ws.Cells.Item(7, 11).Value = 672.4285500000001
$ws.Cells.Item(7, 12).Value = 676.5
$ws.Cells.Item(7, 13).Value = -560.4285500000001
$ws.Cells.Item(7, 14).Value = -900.5
$ws.Cells.Item(19, 8).Value = 300
$ws.Cells.Item(19, 9).Value = 300
$ws.Cells.Item(19, 11).Value = 900
$ws.Cells.Item(19, 13).Value = -726
$ws.Cells.Item(25, 8).Value = 98
$ws.Cells.Item(25, 9).Value = 0
$ws.Cells.Item(25, 10).Value = 98
$ws.Cells.Item(25, 11).Value = 0
$ws.Cells.Item(25, 12).Value = 294
$ws.Cells.Item(25, 13).ClearContents()
$ws.Cells.Item(25, 14).Value = -632
$ws.Cells.Item(30, 8).Value = 98
$ws.Cells.Item(30, 9).Value = 0
$ws.Cells.Item(30, 10).Value = 98
$ws.Cells.Item(30, 11).Value = 0
$ws.Cells.Item(30, 12).Value = 294
$ws.Cells.Item(30, 13).ClearContents()
$ws.Cells.Item(30, 14).Value = -498
$ws.Cells.Item(92, 8).Value = 1072.4286
$ws.Cells.Item(92, 9).Value = 789.6667
$ws.Cells.Item(92, 10).Value = 1284.5
$ws.Cells.Item(92, 11).Value = 2369.0001
$ws.Cells.Item(92, 12).Value = 3853.5
$ws.Cells.Item(92, 13).Value = -1121.0001
$ws.Cells.Item(92, 14).Value = -6349.5
$ws.Cells.Item(107, 8).Value = 12499.667
$ws.Cells.Item(107, 10).Value = 14249.5
$ws.Cells.Item(107, 12).Value = 42748.5
$ws.Cells.Item(107, 14).Value = -46588.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(80, 8).Value = 1405622.4
$ws.Cells.Item(80, 9).Value = 2772912.2
$ws.Cells.Item(80, 11).Value = 2772912.2
$ws.Cells.Item(80, 13).Value = -2771914.2
$ws.Cells.Item(83, 8).Value = 1405622.4
$ws.Cells.Item(83, 9).Value = 2772912.2
$ws.Cells.Item(83, 11).Value = 13864561
$ws.Cells.Item(83, 13).Value = -13859569
$ws.Cells.Item(102, 8).Value = 25009546
$ws.Cells.Item(102, 9).Value = 33342862
$ws.Cells.Item(102, 11).Value = 33342862
$ws.Cells.Item(102, 13).Value = -33341240
$ws.Cells.Item(109, 8).Value = 45000
$ws.Cells.Item(109, 10).Value = 45000
$ws.Cells.Item(109, 12).Value = 45000
$ws.Cells.Item(109, 14).Value = -47080
$ws.Cells.Item(122, 8).Value = 483443.53
$ws.Cells.Item(122, 9).Value = 919017.94
$ws.Cells.Item(122, 10).Value = 8271.454
$ws.Cells.Item(122, 11).Value = 2757053.82
$ws.Cells.Item(122, 12).Value = 24814.362
$ws.Cells.Item(122, 13).Value = -2754603.82
$ws.Cells.Item(122, 14).Value = -29714.362
$ws.Cells.Item(126, 8).Value = 11000
$ws.Cells.Item(126, 9).Value = 3000
$ws.Cells.Item(126, 11).Value = 9000
$ws.Cells.Item(126, 13).Value = -6530
$ws.Cells.Item(132, 8).Value = 101018.234
$ws.Cells.Item(132, 9).Value = 171448.92
$ws.Cells.Item(132, 11).Value = 514346.76
$ws.Cells.Item(132, 13).Value = -511816.76

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(6, 8).Value = 79697
$ws.Cells.Item(6, 10).Value = 79697
$ws.Cells.Item(6, 12).Value = 79697
$ws.Cells.Item(6, 14).Value = -79921
$ws.Cells.Item(22, 8).Value = 996.13336
$ws.Cells.Item(22, 9).Value = 824.75
$ws.Cells.Item(22, 11).Value = 824.75
$ws.Cells.Item(22, 13).Value = -529.75
$ws.Cells.Item(27, 8).Value = 996.13336
$ws.Cells.Item(27, 9).Value = 824.75
$ws.Cells.Item(27, 11).Value = 824.75
$ws.Cells.Item(27, 13).Value = -717.75
$ws.Cells.Item(46, 8).Value = 6699.5654
$ws.Cells.Item(46, 9).Value = 5867
$ws.Cells.Item(46, 10).Value = 6824.45
$ws.Cells.Item(46, 11).Value = 5867
$ws.Cells.Item(46, 12).Value = 6824.45
$ws.Cells.Item(46, 13).Value = -5679
$ws.Cells.Item(46, 14).Value = -7200.45
$ws.Cells.Item(61, 8).Value = 4720.6665
$ws.Cells.Item(61, 9).Value = 4059.6
$ws.Cells.Item(61, 10).Value = 5547
$ws.Cells.Item(61, 11).Value = 4059.6
$ws.Cells.Item(61, 12).Value = 5547
$ws.Cells.Item(61, 13).Value = -3857.6
$ws.Cells.Item(61, 14).Value = -5951
$ws.Cells.Item(100, 8).Value = 7118.5835
$ws.Cells.Item(100, 9).Value = 4680.6
$ws.Cells.Item(100, 10).Value = 8860
$ws.Cells.Item(100, 11).Value = 4680.6
$ws.Cells.Item(100, 12).Value = 8860
$ws.Cells.Item(100, 13).Value = -4139.6
$ws.Cells.Item(100, 14).Value = -9942
$ws.Cells.Item(113, 8).Value = 4720.6665
$ws.Cells.Item(113, 9).Value = 4059.6
$ws.Cells.Item(113, 10).Value = 5547
$ws.Cells.Item(113, 11).Value = 4059.6
$ws.Cells.Item(113, 12).Value = 5547
$ws.Cells.Item(113, 13).Value = -1889.6
$ws.Cells.Item(113, 14).Value = -9887
$ws.Cells.Item(117, 8).Value = 69750
$ws.Cells.Item(117, 10).Value = 90000
$ws.Cells.Item(117, 12).Value = 90000
$ws.Cells.Item(117, 14).Value = -99178
$ws.Cells.Item(131, 8).Value = 66465.664
$ws.Cells.Item(131, 10).Value = 66465.664
$ws.Cells.Item(131, 12).Value = 66465.664
$ws.Cells.Item(131, 14).Value = -76545.664
$ws.Cells.Item(132, 8).Value = 4049.6216
$ws.Cells.Item(132, 9).Value = 2937.68
$ws.Cells.Item(132, 11).Value = 8813.039999999999
$ws.Cells.Item(132, 13).Value = -6283.039999999999
$ws.Cells.Item(136, 8).Value = 3761.5386
$ws.Cells.Item(136, 9).Value = 2322.2222
$ws.Cells.Item(136, 10).Value = 7000
$ws.Cells.Item(136, 11).Value = 6966.6666
$ws.Cells.Item(136, 12).Value = 21000
$ws.Cells.Item(136, 13).Value = -4416.6666
$ws.Cells.Item(136, 14).Value = -26100

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(27, 8).Value = 75624.5
$ws.Cells.Item(27, 10).Value = 75624.5
$ws.Cells.Item(27, 12).Value = 75624.5
$ws.Cells.Item(27, 14).Value = -75762.5
$ws.Cells.Item(34, 8).Value = 50000
$ws.Cells.Item(34, 10).Value = 50000
$ws.Cells.Item(34, 12).Value = 50000
$ws.Cells.Item(34, 14).Value = -50406
$ws.Cells.Item(96, 8).Value = 5279.4546
$ws.Cells.Item(96, 9).Value = 5193.5
$ws.Cells.Item(96, 10).Value = 5328.5713
$ws.Cells.Item(96, 11).Value = 5193.5
$ws.Cells.Item(96, 12).Value = 5328.5713
$ws.Cells.Item(96, 13).Value = -3820.5
$ws.Cells.Item(96, 14).Value = -8074.5713
$ws.Cells.Item(102, 8).Value = 93333.336
$ws.Cells.Item(102, 10).Value = 93333.336
$ws.Cells.Item(102, 12).Value = 93333.336
$ws.Cells.Item(102, 14).Value = -99823.336
$ws.Cells.Item(106, 8).Value = 100000
$ws.Cells.Item(106, 10).Value = 100000
$ws.Cells.Item(106, 12).Value = 100000
$ws.Cells.Item(106, 14).Value = -102524
$ws.Cells.Item(127, 8).Value = 58748.5
$ws.Cells.Item(127, 10).Value = 58748.5
$ws.Cells.Item(127, 12).Value = 58748.5
$ws.Cells.Item(127, 14).Value = -68668.5
$ws.Cells.Item(132, 8).Value = 11114652
$ws.Cells.Item(132, 9).Value = 1244.6666
$ws.Cells.Item(132, 10).Value = 27784764
$ws.Cells.Item(132, 11).Value = 3733.9998
$ws.Cells.Item(132, 12).Value = 83354292
$ws.Cells.Item(132, 13).Value = -1203.9998
$ws.Cells.Item(132, 14).Value = -83359352
$ws.Cells.Item(133, 8).Value = 67568
$ws.Cells.Item(133, 10).Value = 67568
$ws.Cells.Item(133, 12).Value = 67568
$ws.Cells.Item(133, 14).Value = -77688
$ws.Cells.Item(136, 8).Value = 8530.552
$ws.Cells.Item(136, 9).Value = 2827.52
$ws.Cells.Item(136, 11).Value = 8482.56
$ws.Cells.Item(136, 13).Value = -5932.559999999999

